# Updated cryptos list on Sun Jul  7 08:50:55 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns with the latest scrape,
# and fixes rows 41-42 where Filecoin and FirstDigitalUSD had been swapped
# (name, link, price and volume all move back to the correct row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are plain text (coinranking.com scrape),
# including price strings that look numeric, e.g. "57.515.01" or "1.00".
# A bare .Value assignment lets Excel auto-coerce those into real numbers
# (dropping formatting such as trailing zeros), so every write goes
# through this helper: a leading apostrophe forces literal text entry,
# then Style is reset to Normal so no stray "quote prefix" / text format
# is left on the cell (matching the plain, unstyled source cells).
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "57.515.01"
Set-TextValue "E2" "  +1.74%  "
Set-TextValue "D3" "3.010.13"
Set-TextValue "E3" "  +0.33%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "509.91"
Set-TextValue "E5" "  +0.47%  "
Set-TextValue "D6" "139.35"
Set-TextValue "E6" "  +1.78%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "E8" "  +1.28%  "
Set-TextValue "D9" "7.52"
Set-TextValue "E9" "  -0.84%  "
Set-TextValue "E10" "  +1.68%  "
Set-TextValue "D11" "0.365"
Set-TextValue "E11" "  +3.98%  "
Set-TextValue "D12" "3.521.67"
Set-TextValue "E12" "  +0.06%  "
Set-TextValue "E13" "  +0.86%  "
Set-TextValue "D14" "26.47"
Set-TextValue "E14" "  +4.14%  "
Set-TextValue "E15" "  +6.70%  "
Set-TextValue "D16" "57.507.83"
Set-TextValue "E16" "  +1.60%  "
Set-TextValue "D17" "6.21"
Set-TextValue "E17" "  +7.37%  "
Set-TextValue "D18" "3.011.65"
Set-TextValue "E18" "  +0.22%  "
Set-TextValue "D19" "12.82"
Set-TextValue "E19" "  +3.16%  "
Set-TextValue "D20" "7.96"
Set-TextValue "E20" "  +2.04%  "
Set-TextValue "D21" "331.22"
Set-TextValue "E21" "  +1.46%  "
Set-TextValue "D22" "0.997"
Set-TextValue "E22" "  -0.17%  "
Set-TextValue "E23" "  +4.29%  "
Set-TextValue "D24" "64.44"
Set-TextValue "E24" "  +3.29%  "
Set-TextValue "D25" "0.170"
Set-TextValue "E25" "  +0.19%  "
Set-TextValue "E26" "  -0.30%  "
Set-TextValue "D27" "0.0₃0922"
Set-TextValue "E27" "  +1.07%  "
Set-TextValue "D28" "6.80"
Set-TextValue "E28" "  +4.26%  "
Set-TextValue "D29" "7.34"
Set-TextValue "E29" "  +4.77%  "
Set-TextValue "E30" "  +2.16%  "
Set-TextValue "D31" "1.19"
Set-TextValue "E31" "  -5.36%  "
Set-TextValue "D32" "20.60"
Set-TextValue "E32" "  -0.12%  "
Set-TextValue "E33" "  +5.19%  "
Set-TextValue "D34" "153.55"
Set-TextValue "E34" "  -1.37%  "
Set-TextValue "D35" "5.87"
Set-TextValue "E35" "  +4.65%  "
Set-TextValue "D36" "1.28"
Set-TextValue "E36" "  +1.35%  "
Set-TextValue "D37" "0.0683"
Set-TextValue "E37" "  +1.23%  "
Set-TextValue "D38" "24.38"
Set-TextValue "E38" "  +1.07%  "
Set-TextValue "D39" "3.041.04"
Set-TextValue "E39" "  +0.06%  "
Set-TextValue "D40" "37.31"
Set-TextValue "E40" "  +1.80%  "
Set-TextValue "B41" "Filecoin"
Set-TextValue "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "3.85"
Set-TextValue "E41" "  +6.70%  "
Set-TextValue "B42" "FirstDigitalUSD"
Set-TextValue "C42" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.11%  "
Set-TextValue "E43" "  -0.12%  "
Set-TextValue "D44" "2.247.73"
Set-TextValue "E44" "  -0.89%  "
Set-TextValue "E45" "  +0.58%  "
Set-TextValue "E46" "  -1.03%  "
Set-TextValue "D47" "6.02"
Set-TextValue "E47" "  +4.82%  "
Set-TextValue "E48" "  +1.94%  "
Set-TextValue "D49" "19.39"
Set-TextValue "E49" "  +2.06%  "
Set-TextValue "E50" "  -6.54%  "
Set-TextValue "E51" "  +2.66%  "
